# Generate Report for Handoff
# Updates the localization-status workbook from a "handed back" snapshot to a
# fresh "ready for handoff" snapshot: new source/handoff file identifiers,
# new handoff timestamps, status reset to "Ready for handoff", the handback
# (target/handback-file/handback-datetime) columns cleared out on the
# per-language detail sheets, and both rows now point at the same handoff
# package.

$wb = $excel.ActiveWorkbook

$oldMdA = "8af56c6e-9cba-47e7-b99d-f0ce945b5201.md"
$oldMdB = "a8e99754-0d9f-4d5c-9714-60afd90a4c49.md"
$newMdA = "34ded686-4006-40a0-a24a-57ef94237596.md"
$newMdB = "ffff11cb6873-c225-4731-8417-48de30a3441c.md"

$newXlfZh = "34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf"
$newXlfDe = "34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf"

$newStatus = "Ready for handoff"
$newHandoffDateOverview = "2016-03-24 08:49:21"
$newHandoffDatetimeZh = "2016-03-24 08:49:16"
$newHandoffDatetimeDe = "2016-03-24 08:49:21"
$resetHandbackDatetime = "0001-01-01 00:00:00"

function Set-HyperlinkDisplay($ws, [string]$addr, [string]$text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
            return
        }
    }
}

function Remove-HyperlinkAt($ws, [string]$addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            return
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet: one row per source file, rolled-up status/date columns.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdA
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = $newHandoffDateOverview

$wsOverview.Range("A3").Value = $newMdB
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $newHandoffDateOverview

Set-HyperlinkDisplay $wsOverview '$A$2' $newMdA
Set-HyperlinkDisplay $wsOverview '$A$3' $newMdB

# ---------------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdA
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $newHandoffDatetimeZh
$wsZh.Range("H2").Value = $resetHandbackDatetime

$wsZh.Range("A3").Value = $newMdB
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $newHandoffDatetimeZh
$wsZh.Range("H3").Value = $resetHandbackDatetime

Set-HyperlinkDisplay $wsZh '$A$2' $newMdA
Set-HyperlinkDisplay $wsZh '$D$2' $newXlfZh
Set-HyperlinkDisplay $wsZh '$A$3' $newMdB
Set-HyperlinkDisplay $wsZh '$D$3' $newXlfZh

Remove-HyperlinkAt $wsZh '$F$2'
Remove-HyperlinkAt $wsZh '$G$2'
Remove-HyperlinkAt $wsZh '$F$3'
Remove-HyperlinkAt $wsZh '$G$3'

$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()

# ---------------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdA
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $newHandoffDatetimeDe
$wsDe.Range("H2").Value = $resetHandbackDatetime

$wsDe.Range("A3").Value = $newMdB
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $newHandoffDatetimeDe
$wsDe.Range("H3").Value = $resetHandbackDatetime

Set-HyperlinkDisplay $wsDe '$A$2' $newMdA
Set-HyperlinkDisplay $wsDe '$D$2' $newXlfDe
Set-HyperlinkDisplay $wsDe '$A$3' $newMdB
Set-HyperlinkDisplay $wsDe '$D$3' $newXlfDe

Remove-HyperlinkAt $wsDe '$F$2'
Remove-HyperlinkAt $wsDe '$G$2'
Remove-HyperlinkAt $wsDe '$F$3'
Remove-HyperlinkAt $wsDe '$G$3'

$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()
